# Style all data added
# Prefix each category label in column A (rows 2-11) with a letter A-J,
# reflecting newly-added styled data for each series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefixes = @("A","B","C","D","E","F","G","H","I","J")

for ($i = 0; $i -lt $prefixes.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value2
    $cell.Value = "$($prefixes[$i])-$current"
}

# Update the active selection to A11, matching the final selected cell.
$ws.Range("A11").Select()
